$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), copying the formatting
# (bold font + border + alignment) already used by the other header
# cells (e.g. H1) so the new headers match the existing style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the data values for columns I (I0) and J (IF) for rows 2-33.
$iValues = @(8,8,8,8,9,8,6,9,8,9,7,9,8,9,8,7,5,7,8,5,9,3,5,7,6,8,7,7,9,7,7,3)
$jValues = @(8,8,8,8,9,8,6,9,8,9,8,9,9,9,8,8,5,7,9,5,9,4,6,7,6,8,8,7,9,7,7,3)

for ($n = 0; $n -lt 32; $n++) {
    $r = $n + 2
    $ws.Cells.Item($r, 9).Value = $iValues[$n]
    $ws.Cells.Item($r, 10).Value = $jValues[$n]
}

Write-Output "I0/IF columns added"
